$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Append the new sentences after "Vi tar pause til å spise lunsj. "
#    in the Wednesday (18.1.2023) paragraph, as eight separate runs
#    (matching how the diff records them).
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*Vi tar pause til å spise lunsj.*") {
        $target = $para
    }
}

$r = $target.Range
$r.MoveEnd(1, -1)
$r.Collapse(0)

$r.InsertAfter("Etter lunsj så kom Viet-Uy, Besnik, Ilyas, Mathias og Fergus")
$r.Collapse(0)

$r.InsertAfter(".")
$r.Collapse(0)

$r.InsertAfter(" Viet-Uy")
$r.Collapse(0)

$r.InsertAfter(" måtte gå til å ha en time.")
$r.Collapse(0)

$r.InsertAfter(" ")
$r.Collapse(0)

$r.InsertAfter("Ilyas og Besnik jobber med ")
$r.Collapse(0)

$r.InsertAfter("å utvikle nettsiden mer")
$r.Collapse(0)

$r.InsertAfter(", alle andre fortsetter med hva de jobbet med. ")
$r.Collapse(0)

# Stamp the whole paragraph with the nb-NO language so the freshly
# inserted runs carry <w:lang w:val="nb-NO"/> just like their
# neighbours (a no-op for the pre-existing runs, which already are
# nb-NO).
$target.Range.LanguageID = "nb-NO"

# ------------------------------------------------------------------
# 2) Drop the stray <w:lastRenderedPageBreak/> in front of the second
#    "[placeholder]" run (Friday 20.1.2023 section).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*[placeholder]*") {
        $pr = $para.Range
        $pr.Find.Execute("[placeholder]", $true, $false, $false, $false, $false, $true, 1, $false, "[placeholder]", 2) | Out-Null
    }
}
